$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.100996534501604
$ws.Range("E5").Value = 0.9252513742688047
$ws.Range("H5").Value = 0.03674991256134175
$ws.Range("K5").Value = 2.022620531650275
$ws.Range("N5").Value = 1.5435538883144433
$ws.Range("Q5").Value = 3.217584363748626
$ws.Range("T5").Value = 1.276640880978019
$ws.Range("W5").Value = 7.310420465267463
$ws.Range("Z5").Value = 3.958338677803065
$ws.Range("AC5").Value = 1.435186030802653
$ws.Range("AF5").Value = 8.651470290175885
$ws.Range("AI5").Value = 6.065336416975523
$ws.Range("AL5").Value = 2.025690748957605
$ws.Range("AO5").Value = 1.546103239782389
$ws.Range("AR5").Value = 11.305240510113663

$ws = $wb.Worksheets.Item(2)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 4.888873861980685
$ws.Range("E5").Value = 0.9185060081482374
$ws.Range("H5").Value = 0.07611365453726998
$ws.Range("K5").Value = 3.6079480655280682
$ws.Range("N5").Value = 1.3327697571534383
$ws.Range("Q5").Value = 18.314914732994406
$ws.Range("T5").Value = 4.647029500955299
$ws.Range("W5").Value = 61.648661267478985
$ws.Range("Z5").Value = 19.192373444759014
$ws.Range("AC5").Value = 0.9418775679608319
$ws.Range("AF5").Value = 12.728459148555919
$ws.Range("AI5").Value = 3.7688273712551523
$ws.Range("AL5").Value = 3.5719138219381095
$ws.Range("AO5").Value = 1.331607527895039
$ws.Range("AR5").Value = 9.864629574212351

$ws = $wb.Worksheets.Item(3)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 5.307886017380093
$ws.Range("E5").Value = 0.8382896808494182
$ws.Range("H5").Value = 0.14798131495278358
$ws.Range("K5").Value = 1.310569779336032
$ws.Range("N5").Value = 2.497587087736609
$ws.Range("Q5").Value = 17.994745137005538
$ws.Range("T5").Value = 6.096536096521219
$ws.Range("W5").Value = 14.227817447783437
$ws.Range("Z5").Value = 20.949121620098417
$ws.Range("AC5").Value = 0.9587516268514867
$ws.Range("AF5").Value = 3.7284382659756954
$ws.Range("AI5").Value = 5.4545855893714785
$ws.Range("AL5").Value = 1.3434662272539215
$ws.Range("AO5").Value = 2.523244536051863
$ws.Range("AR5").Value = 6.2828715262505845

$ws = $wb.Worksheets.Item(4)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.9391982816475433
$ws.Range("E5").Value = 1.357056631934884
$ws.Range("H5").Value = 0.017455406595220976
$ws.Range("K5").Value = 4.062573276659963
$ws.Range("N5").Value = 0.31733747346084773
$ws.Range("Q5").Value = 4.65228064271653
$ws.Range("T5").Value = 0.8653972805485179
$ws.Range("W5").Value = 19.224383539444602
$ws.Range("Z5").Value = 1.2578035787494988
$ws.Range("AC5").Value = 1.1785833595356252
$ws.Range("AF5").Value = 21.452901543195754
$ws.Range("AI5").Value = 1.8149147875522795
$ws.Range("AL5").Value = 4.059574799504038
$ws.Range("AO5").Value = 0.3170688175179208
$ws.Range("AR5").Value = 6.298669576644133

$ws = $wb.Worksheets.Item(5)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.4413483752713354
$ws.Range("E5").Value = 1.1460439912611025
$ws.Range("H5").Value = 0.012267069208337129
$ws.Range("K5").Value = 3.305532869948481
$ws.Range("N5").Value = 0.6109842493858042
$ws.Range("Q5").Value = 1.6492372865760085
$ws.Range("T5").Value = 0.43654673739990746
$ws.Range("W5").Value = 5.475535116337326
$ws.Range("Z5").Value = 1.237532109062851
$ws.Range("AC5").Value = 1.1590204238854565
$ws.Range("AF5").Value = 11.277081911577866
$ws.Range("AI5").Value = 2.817959245252643
$ws.Range("AL5").Value = 3.3071136264460885
$ws.Range("AO5").Value = 0.6114745493827558
$ws.Range("AR5").Value = 9.19244904166535

$ws = $wb.Worksheets.Item(6)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 4.887848406532721
$ws.Range("E5").Value = 1.0098936820016218
$ws.Range("H5").Value = 0.07985274760400614
$ws.Range("K5").Value = 7.063697007392537
$ws.Range("N5").Value = 2.034875456516579
$ws.Range("Q5").Value = 22.051338859385854
$ws.Range("T5").Value = 5.03450716263627
$ws.Range("W5").Value = 105.7069295837905
$ws.Range("Z5").Value = 28.215313284379913
$ws.Range("AC5").Value = 1.220679843682869
$ws.Range("AF5").Value = 30.139900523461094
$ws.Range("AI5").Value = 6.028414426499104
$ws.Range("AL5").Value = 7.056236982776979
$ws.Range("AO5").Value = 2.038329952925962
$ws.Range("AR5").Value = 21.88827727219915

$ws = $wb.Worksheets.Item(7)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.31785982940659036
$ws.Range("E5").Value = 1.0856945208808453
$ws.Range("H5").Value = 0.006384554910962774
$ws.Range("K5").Value = 3.5319821164977054
$ws.Range("N5").Value = 0.6143958688116288
$ws.Range("Q5").Value = 1.2574691018038093
$ws.Range("T5").Value = 0.1893463129014461
$ws.Range("W5").Value = 4.583081122684927
$ws.Range("Z5").Value = 0.8237596982663389
$ws.Range("AC5").Value = 1.3932112661257934
$ws.Range("AF5").Value = 24.992122386190896
$ws.Range("AI5").Value = 4.875920560724343
$ws.Range("AL5").Value = 3.532244826646799
$ws.Range("AO5").Value = 0.613786803608615
$ws.Range("AR5").Value = 27.58089957520527

$ws = $wb.Worksheets.Item(8)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.8386692362453372
$ws.Range("E5").Value = 1.4750939701833212
$ws.Range("H5").Value = 0.019125834124636466
$ws.Range("K5").Value = 4.499084244048931
$ws.Range("N5").Value = 0.6536122775618325
$ws.Range("Q5").Value = 10.221638988585704
$ws.Range("T5").Value = 1.7708632956253358
$ws.Range("W5").Value = 37.43665188789094
$ws.Range("Z5").Value = 6.888482111948018
$ws.Range("AC5").Value = 1.4496819035045538
$ws.Range("AF5").Value = 25.443320522046324
$ws.Range("AI5").Value = 3.687616260752303
$ws.Range("AL5").Value = 4.475323719259719
$ws.Range("AO5").Value = 0.6497162088230137
$ws.Range("AR5").Value = 9.554539929195615

$ws = $wb.Worksheets.Item(9)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.5585429628976407
$ws.Range("E5").Value = 1.196348976447853
$ws.Range("H5").Value = 0.014571713361833431
$ws.Range("K5").Value = 4.830991999055013
$ws.Range("N5").Value = 2.2508825618988673
$ws.Range("Q5").Value = 1.4914592635149895
$ws.Range("T5").Value = 0.4011394454699258
$ws.Range("W5").Value = 13.091799168629793
$ws.Range("Z5").Value = 11.518684119529617
$ws.Range("AC5").Value = 0.9494370298081568
$ws.Range("AF5").Value = 61.050731876527976
$ws.Range("AI5").Value = 55.514936810659655
$ws.Range("AL5").Value = 4.798659683322402
$ws.Range("AO5").Value = 2.2204159798148164
$ws.Range("AR5").Value = 204.32935788696545

$ws = $wb.Worksheets.Item(10)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.17060424980598507
$ws.Range("E5").Value = 0.833734662982283
$ws.Range("H5").Value = 0.00266304614972239
$ws.Range("K5").Value = 1.897484167180915
$ws.Range("N5").Value = 0.017639066039341886
$ws.Range("Q5").Value = 0.9900985522397997
$ws.Range("T5").Value = 0.18156923211239628
$ws.Range("W5").Value = 2.292415367705119
$ws.Range("Z5").Value = 0.19155259682525203
$ws.Range("AC5").Value = 1.049279705579732
$ws.Range("AF5").Value = 8.2878600219588
$ws.Range("AI5").Value = 1.1151629215969892
$ws.Range("AL5").Value = 1.8982001914642248
$ws.Range("AO5").Value = 0.017628022014683078
$ws.Range("AR5").Value = 1.8966327963876983

$ws = $wb.Worksheets.Item(11)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.4285852477327999
$ws.Range("E5").Value = 1.1592346584283886
$ws.Range("H5").Value = 0.006117605554676481
$ws.Range("K5").Value = 4.599061586492872
$ws.Range("N5").Value = 1.289878090575323
$ws.Range("Q5").Value = 2.4684143479180283
$ws.Range("T5").Value = 0.409897929326988
$ws.Range("W5").Value = 14.011633755990408
$ws.Range("Z5").Value = 6.917295613575523
$ws.Range("AC5").Value = 0.8813048802602795
$ws.Range("AF5").Value = 33.95230965473998
$ws.Range("AI5").Value = 14.899430842038095
$ws.Range("AL5").Value = 4.599620307055454
$ws.Range("AO5").Value = 1.2910407119278575
$ws.Range("AR5").Value = 84.81120603692578

$ws = $wb.Worksheets.Item(12)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.2908855775859471
$ws.Range("E5").Value = 1.4221735287880926
$ws.Range("H5").Value = 0.005043524273514312
$ws.Range("K5").Value = 2.8922437684444335
$ws.Range("N5").Value = 1.5803709201301404
$ws.Range("Q5").Value = 1.9929829546965563
$ws.Range("T5").Value = 0.2704956872828417
$ws.Range("W5").Value = 7.718497221657607
$ws.Range("Z5").Value = 3.3622866089883323
$ws.Range("AC5").Value = 1.6233509392818164
$ws.Range("AF5").Value = 36.26470478265078
$ws.Range("AI5").Value = 17.000701249723974
$ws.Range("AL5").Value = 2.8896651568185323
$ws.Range("AO5").Value = 1.5783120111480897
$ws.Range("AR5").Value = 66.02600653852342
